# Updates cryptos list: refresh prices / 1h volume deltas; swap the
# Bittensor/OKB ranking rows (#41/#42) to match the new ordering.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '68.586.70'
$ws.Range("E2").Value = '  +1.80%  '

# Row 3
$ws.Range("D3").Value = '3.264.56'
$ws.Range("E3").Value = '  +0.56%  '

# Row 4
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.96'
$ws.Range("E5").Value = '  +1.04%  '

# Row 6
$ws.Range("E6").Value = '  -0.04%  '

# Row 7
$ws.Range("E7").Value = '  +0.00%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.132'
$ws.Range("E9").Value = '  +0.67%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.65'
$ws.Range("E10").Value = '  -1.39%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.423'
$ws.Range("E11").Value = '  +2.11%  '

# Row 12
$ws.Range("D12").Value = '3.834.29'
$ws.Range("E12").Value = '  +0.61%  '

# Row 13
$ws.Range("E13").Value = '  +0.14%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.51'
$ws.Range("E14").Value = '  +0.21%  '

# Row 15
$ws.Range("D15").Value = '68.557.55'
$ws.Range("E15").Value = '  +1.80%  '

# Row 16
$ws.Range("E16").Value = '  +2.26%  '

# Row 17
$ws.Range("D17").Value = '3.209.97'
$ws.Range("E17").Value = '  -1.12%  '

# Row 18
$ws.Range("E18").Value = '  -0.38%  '

# Row 19
$ws.Range("E19").Value = '  +0.42%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '393.56'
$ws.Range("E20").Value = '  +4.64%  '

# Row 21
$ws.Range("E21").Value = '  +1.11%  '

# Row 22
$ws.Range("E22").Value = '  +1.13%  '

# Row 23
$ws.Range("E23").Value = '  -0.17%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.515'
$ws.Range("E24").Value = '  +0.93%  '

# Row 25
$ws.Range("E25").Value = '  +0.33%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.189'
$ws.Range("E26").Value = '  +4.33%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.59'
$ws.Range("E27").Value = '  +0.32%  '

# Row 28
$ws.Range("E28").Value = '  -0.52%  '

# Row 29
$ws.Range("E29").Value = '  +0.29%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.68'
$ws.Range("E30").Value = '  -1.57%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.91'
$ws.Range("E31").Value = '  +1.16%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.12'

# Row 33
$ws.Range("E33").Value = '  +0.82%  '

# Row 34
$ws.Range("E34").Value = '  +0.03%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '164.13'
$ws.Range("E35").Value = '  +0.08%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.50'
$ws.Range("E36").Value = '  +0.34%  '

# Row 37
$ws.Range("E37").Value = '  +3.38%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.828'
$ws.Range("E38").Value = '  -2.24%  '

# Row 39
$ws.Range("E39").Value = '  -0.37%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '26.30'
$ws.Range("E40").Value = '  -1.25%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.55'
$ws.Range("E41").Value = '  -3.97%  '

# Row 42
$ws.Range("E42").Value = '  -4.15%  '

# Row 43
$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '347.97'
$ws.Range("E43").Value = '  -2.93%  '

# Row 44
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.28'
$ws.Range("E44").Value = '  +1.13%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0687'
$ws.Range("E45").Value = '  +1.04%  '

# Row 46
$ws.Range("D46").Value = '2.610.47'
$ws.Range("E46").Value = '  -3.87%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.59'
$ws.Range("E47").Value = '  -3.66%  '

# Row 48
$ws.Range("E48").Value = '  +0.82%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.31'
$ws.Range("E49").Value = '  +2.77%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '31.55'
$ws.Range("E50").Value = '  +1.10%  '

# Row 51
$ws.Range("E51").Value = '  -0.32%  '
